$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-10-19"

# Update the header text for the current-month column (B1)
$ws.Range("B1").Value = "October 2021 (through October 19)"

# Garfield Park (row 2)
$ws.Range("B2").Value = 16
$ws.Range("L2").Value = 12

# Auburn Gresham (row 6)
$ws.Range("B6").Value = 5

# Little Italy, UIC (row 8) - new value
$ws.Range("B8").Value = 1

# South Shore (row 13)
$ws.Range("B13").Value = 6

# United Center (row 15) - new value
$ws.Range("AF15").Value = 1

# Grand Boulevard (row 22)
$ws.Range("B22").Value = 3

# Chinatown (row 24)
$ws.Range("B24").Value = 3

# Avondale (row 36)
$ws.Range("B36").Value = 2

# Englewood (row 38)
$ws.Range("AP38").Value = 3

# Washington Heights (row 41)
$ws.Range("L41").Value = 5

# Washington Park (row 42) - new value
$ws.Range("B42").Value = 1

# Bridgeport (row 48)
$ws.Range("L48").Value = 3
$ws.Range("BJ48").Value = 1

# Gage Park (row 72)
$ws.Range("L72").Value = 2

# Near South Side (row 87) - new value
$ws.Range("B87").Value = 1

# Woodlawn (row 98) - new value
$ws.Range("AZ98").Value = 1
